$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new survey rows (77-117) for participants 4iaixd9p, 5virtmb4, gtot06qe
$data = @(
    ,@('4iaixd9p', 'Training phase', 1, '[''Purple'', ''Orange'', ''Green'']', '[[''Red'', ''''], [''Blue'', ''''], [''Blue'', ''Blue'']]')
    ,@('4iaixd9p', 'Training phase', 2, '[''Green'', ''Green'', ''Orange'']', '[[''Yellow'', ''''], [''Blue'', ''''], [''Yellow'', ''Red'']]')
    ,@('4iaixd9p', 'Training phase', 3, '[''Purple'', ''Green'', ''Purple'']', '[[''Red'', ''Blue''], [''Blue'', ''Yellow''], [''Red'', ''Blue'']]')
    ,@('4iaixd9p', 'Training phase', 4, '[''Purple'', ''Orange'', ''Orange'']', '[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''Yellow'']]')
    ,@('4iaixd9p', 'Training phase', 5, '[''Green'', ''Orange'', ''Purple'']', '[[''Red'', ''Red''], [''Blue'', ''''], [''Red'', ''Yellow'']]')
    ,@('4iaixd9p', 'Test 1', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Yellow'', ''Blue''], [''Yellow'', ''Yellow''], [''Red'', ''Blue''], [''Red'', ''Red''], [''Red'', ''Yellow''], [''Blue'', ''Blue'']]')
    ,@('4iaixd9p', 'Exploration phase', 1, '[''Green'', ''Purple'', ''Blue'']', '[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('4iaixd9p', 'Exploration phase', 2, '[''Green'', ''Purple'', ''Blue'']', '[[''Red'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('4iaixd9p', 'Exploration phase', 3, '[''Green'', ''Purple'', ''Yellow'']', '[[''Blue'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('4iaixd9p', 'Exploration phase', 4, '[''Green'', ''Purple'', ''Yellow'']', '[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('4iaixd9p', 'Exploration phase', 5, '[''Green'', ''Purple'', ''Yellow'']', '[[''Blue'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('4iaixd9p', 'Test 2', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Blue'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Yellow'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('4iaixd9p', 'Test 2', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Red'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('4iaixd9p', 'Test 2', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Blue'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Blue'', '''']]')
    ,@('4iaixd9p', 'Test 2', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Yellow'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Red'', ''Red'']]')
    ,@('4iaixd9p', 'Test 2', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Blue'', ''''], [''Red'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('4iaixd9p', 'Test 2', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Blue'', ''''], [''Red'', ''''], [''Red'', ''''], [''Yellow'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('4iaixd9p', 'Test 2', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Yellow'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Yellow'', ''''], [''Red'', ''''], [''Blue'', '''']]')
    ,@('5virtmb4', 'Training phase', 1, '[''Purple'', ''Orange'', ''Green'']', '[[''Yellow'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('5virtmb4', 'Training phase', 2, '[''Green'', ''Green'', ''Orange'']', '[[''Red'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('5virtmb4', 'Training phase', 3, '[''Purple'', ''Green'', ''Purple'']', '[[''Red'', ''''], [''Yellow'', ''''], [''Blue'', '''']]')
    ,@('5virtmb4', 'Training phase', 4, '[''Purple'', ''Orange'', ''Orange'']', '[[''Yellow'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('5virtmb4', 'Training phase', 5, '[''Green'', ''Orange'', ''Purple'']', '[[''Blue'', ''''], [''Yellow'', ''''], [''Blue'', '''']]')
    ,@('5virtmb4', 'Test 1', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Yellow'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('5virtmb4', 'Test 1', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Blue'', ''''], [''Red'', ''''], [''Red'', ''''], [''Yellow'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('5virtmb4', 'Test 1', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Blue'', ''''], [''Red'', ''''], [''Red'', ''''], [''Yellow'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('5virtmb4', 'Test 1', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Red'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Red'', '''']]')
    ,@('5virtmb4', 'Test 1', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('5virtmb4', 'Exploration phase', 1, '[''Green'', ''Purple'', ''Blue'']', '[[''Red'', ''''], [''Red'', ''''], [''Red'', '''']]')
    ,@('5virtmb4', 'Exploration phase', 2, '[''Green'', ''Purple'', ''Blue'']', '[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('5virtmb4', 'Exploration phase', 3, '[''Green'', ''Purple'', ''Yellow'']', '[[''Red'', ''''], [''Yellow'', ''''], [''Blue'', '''']]')
    ,@('5virtmb4', 'Exploration phase', 4, '[''Green'', ''Purple'', ''Yellow'']', '[[''Blue'', ''''], [''Red'', ''''], [''Blue'', '''']]')
    ,@('5virtmb4', 'Exploration phase', 5, '[''Green'', ''Purple'', ''Yellow'']', '[[''Blue'', ''''], [''Yellow'', ''''], [''Blue'', '''']]')
    ,@('5virtmb4', 'Test 2', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Blue'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('5virtmb4', 'Test 2', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Yellow'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('gtot06qe', 'Training phase', 1, '[''Purple'', ''Orange'', ''Green'']', '[[''Blue'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('gtot06qe', 'Training phase', 2, '[''Green'', ''Green'', ''Orange'']', '[[''Yellow'', ''''], [''Red'', ''''], [''Blue'', '''']]')
    ,@('gtot06qe', 'Training phase', 3, '[''Purple'', ''Green'', ''Purple'']', '[[''Red'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('gtot06qe', 'Training phase', 4, '[''Purple'', ''Orange'', ''Orange'']', '[[''Blue'', ''''], [''Red'', ''''], [''Yellow'', '''']]')
    ,@('gtot06qe', 'Training phase', 5, '[''Green'', ''Orange'', ''Purple'']', '[[''Yellow'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('gtot06qe', 'Test 1', 1, '[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']', '[[''Yellow'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Red'', '''']]')
)

$startRow = 77
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Widen column E (nutrients) from 110 to 114 raw OOXML width units
$ws.Columns.Item(5).ColumnWidth = 113.16666666666667
